# Actualización SmartScore desde Streamlit (Salvador Vidal)
# Adds three new header columns (Nombre Completo / Edad / Género),
# converts the Top-3 SmartScore cells on row 3 from text to real numbers,
# and appends a brand-new survey response row (row 4) for Salvador Vidal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New header cells AE1:AG1 — copy the look of the existing header
#    (bold font, border, centered/top alignment) then set their text.
# ---------------------------------------------------------------------
$ws.Range("AD1").Copy()
$ws.Range("AE1:AG1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AE1").Value = "Nombre Completo"
$ws.Range("AF1").Value = "Edad"
$ws.Range("AG1").Value = "Género"

# ---------------------------------------------------------------------
# 2) Row 3: the Top-3 SmartScore cells were stored as text — convert
#    them to actual numeric values.
# ---------------------------------------------------------------------
$ws.Range("E3").Value = 0.533
$ws.Range("H3").Value = 0.475
$ws.Range("K3").Value = 0.426
$ws.Range("N3").Value = 0.703
$ws.Range("Q3").Value = 0.639
$ws.Range("T3").Value = 0.552
$ws.Range("W3").Value = 0.698
$ws.Range("Z3").Value = 0.601
$ws.Range("AC3").Value = 0.579

# ---------------------------------------------------------------------
# 3) Append row 4 — new Streamlit submission from Salvador Vidal.
#    Helper: write a numeric-looking string while keeping it as TEXT
#    (matches the source export, which keeps SmartScore values as text
#    for this respondent) without leaving a stray number-format style.
# ---------------------------------------------------------------------
function Set-TextValue($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("A4").Value = "Salvador Vidal"
$ws.Range("B4").Value = "2025-10-28 04:13:10"

$json = @"
{
  "portion": 0.8,
  "diet": 0.42857142857142855,
  "salt": 0.4,
  "fat": 0.8,
  "natural": 0.8,
  "convenience": 0.2,
  "price": 0.6
}
"@
$ws.Range("C4").Value = $json

$ws.Range("D4").Value = "Nongshim Neoguri Spicy Seafood"
Set-TextValue "E4" "0.627"
$ws.Range("F4").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

$ws.Range("G4").Value = "Nissin Chow Mein Teriyaki Beef"
Set-TextValue "H4" "0.469"
$ws.Range("I4").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

$ws.Range("J4").Value = "Maruchan Ramen Sabor Pollo"
Set-TextValue "K4" "0.458"
$ws.Range("L4").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"

$ws.Range("M4").Value = "Kraft Macaroni & Cheese Dinner"
Set-TextValue "N4" "0.696"
$ws.Range("O4").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

$ws.Range("P4").Value = "Annie’s Shells & White Cheddar"
Set-TextValue "Q4" "0.640"
$ws.Range("R4").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"

$ws.Range("S4").Value = "Amy’s Macaroni & Cheese (frozen)"
Set-TextValue "T4" "0.577"
$ws.Range("U4").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"

$ws.Range("V4").Value = "Wild Planet Wild Tuna Pasta Salad"
Set-TextValue "W4" "0.679"
$ws.Range("X4").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

$ws.Range("Y4").Value = "Kitchens of India Variety Pack"
Set-TextValue "Z4" "0.545"
$ws.Range("AA4").Value = "Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad"

$ws.Range("AB4").Value = "StarKist Chicken Creations (Chicken Salad)"
Set-TextValue "AC4" "0.516"
$ws.Range("AD4").Value = "Portátil, saludable, fácil, buena textura, sabor suave"

$ws.Range("AE4").Value = "Salvador Vidal"
$ws.Range("AF4").Value = 22
$ws.Range("AG4").Value = "Masculino"

# Setting the multi-line JSON text in C4 makes Excel auto-apply a custom
# row height; AutoFit() resets row 4 back to the sheet's default height
# (matching rows 1-3, which have no explicit row height either).
$ws.Rows.Item(4).AutoFit()
